$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the existing header style
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting of the existing header cell (H1) onto the new header
# cells so they reuse the same bold/centered/bordered style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data columns I and J for rows 2-4
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 6
